$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Users")

$ws.Range("A3").Value = "moshe1"
$ws.Range("B3").Value = "Mm1122334!"
$ws.Range("C3").Value = 206676850
$ws.Range("D3").Value = "moses@gmail.com"
$ws.Range("E3").Value = "boy"
$ws.Range("F3").Value = 0
